# Commit: "Added policyForm, Year built, older Reneovated Home"
#
# Sheet1 gains three new trailing columns (N=YearBuilt, O=ReneovatedHome,
# P=PolicyForm); the existing address/city/state/zipcode/occupancy/
# fire-hydrant-rate sample record is swapped for a different sample
# record; Sheet2 gains a small 1-row/4-column reference table holding the
# "old" address that used to live in Sheet1.
#
# The writes below are ordered to match the order the strings were first
# typed (so new entries land in xl/sharedStrings.xml in the same order
# as the target file).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- New header cells on Sheet1 (N1:P1) ---
$ws1.Range("N1").Value = "YearBuilt"
$ws1.Range("O1").Value = "ReneovatedHome"

# --- New little lookup table on Sheet2 (the "old" address) ---
$ws2.Range("A2").Value = "195 Painted Desert Ln"
$ws2.Range("B2").Value = "Buda"
$ws2.Range("C2").Value = "TX"
$ws2.Range("D2").Value = 78610

# --- Sheet1 row 2 data updates ---
$ws1.Range("M2").Value = "YES"

$ws1.Range("P1").Value = "PolicyForm"
$ws1.Range("P2").Value = "DP-3"

$ws1.Range("L2").Value = "Seasonal"

$ws1.Range("G2").Value = "312 Saddle Wood Dr"
$ws1.Range("H2").Value = "Canton"
$ws1.Range("I2").Value = "WI"

$ws1.Range("O2").Value = "YES"

# DOB/Zipcode-style numeric cell (keeps its original "0;[Red]0" look),
# now repurposed to hold a date serial.
$ws1.Range("J2").Value = 30114

# YearBuilt value - copy J2's number style (0;[Red]0) across so it isn't
# re-interpreted as text (column N inherits the sheet's text style).
$ws1.Range("J2").Copy()
$ws1.Range("N2").PasteSpecial(-4122)
$ws1.Range("N2").Value = 2000

# ReneovatedHome (O2) needs its own "text" style (numFmt 49, no center
# alignment) - start from the default "Normal" style then apply the text
# number format, which mints exactly that new cellXf.
$ws1.Range("O2").Style = "Normal"
$ws1.Range("O2").NumberFormat = "@"
$ws1.Range("O2").Value = "YES"

# Sheet2's new cells should look like the rest of the workbook: text
# cells centered like column A1's style, the zip like J2's number style.
$ws1.Range("A1").Copy()
$ws2.Range("A2:C2").PasteSpecial(-4122)
$ws1.Range("J2").Copy()
$ws2.Range("D2").PasteSpecial(-4122)

# Column widths for the two new "bestFit" columns on Sheet1.
$ws1.Columns.Item(15).ColumnWidth = 16.42
$ws1.Columns.Item(16).ColumnWidth = 9.92

# Sheet2's first column was widened to fit the new address text.
$ws2.Columns.Item(1).ColumnWidth = 19.5

# Leave the cursor/selection roughly where the author left it: full used
# range selected on Sheet1 (the active tab), a single highlighted row on
# Sheet2.
$ws2.Range("A2:D2").Select()
$ws1.Activate()
$ws1.Cells.Select()
